$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new cell values in column C next to the existing item rows (5, 6, 7)
# so they line up with the "Status" / "DONE" header cells already present
# at C1 / C4. These new cells reuse the "DONE" text and introduce two new
# shared strings ("most" and "just design").
$ws.Range("C7").Value = "just design"
$ws.Range("C6").Value = "DONE"
$ws.Range("C5").Value = "most"
